# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Coliflor"
# right before the existing row 179. This shifts all subsequent rows
# (old rows 179-205) down by one (new rows 180-206), and the new row 179
# is populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 179; existing rows 179+ shift down to 180+
$ws.Rows.Item(179).Insert()

# Copy the full row that is now 180 (the former row 179) into the newly
# inserted row 179 so the new row inherits matching cell formatting
# (in particular the date style used in column D).
$ws.Range("A180:R180").Copy()
$ws.Range("A179").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Now overwrite with the new record's values
$ws.Cells.Item(179, 1).Value = 4
$ws.Cells.Item(179, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(179, 3).Value = "Los Lagos"
$ws.Cells.Item(179, 4).Value = 44491
$ws.Cells.Item(179, 5).Value = 10
$ws.Cells.Item(179, 6).Value = 100112008
$ws.Cells.Item(179, 7).Value = "Coliflor"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 1400
$ws.Cells.Item(179, 11).Value = 1200
$ws.Cells.Item(179, 12).Value = 1300
$ws.Cells.Item(179, 13).Value = 1250
$ws.Cells.Item(179, 14).Value = "`$/unidad"
$ws.Cells.Item(179, 15).Value = "Región Metropolitana"
$ws.Cells.Item(179, 16).Value = 1250
$ws.Cells.Item(179, 17).Value = 1
$ws.Cells.Item(179, 18).Value = "Hortaliza"
